$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (ten -> nine data rows after the TPM re-run)
$ws.Rows.Item(10).Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ptn"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.327816333333333
$ws.Range("H2").Value = 6.983449
$ws.Range("I2").Value = 0.02128501190197005
$ws.Range("J2").Value = 0.02128501190197004
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05372733333333333
$ws.Range("N2").Value = 0.161182
$ws.Range("O2").Value = 0.1072370469527173
$ws.Range("P2").Value = 0.1072370469527173
$ws.Range("Q2").Value = 0.1250673640797778
$ws.Range("R2").Value = 1.125606276718
$ws.Range("S2").Value = 0.002282541820720708
$ws.Range("T2").Value = 0.002282541820720708

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ptn"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.327816333333333
$ws.Range("H3").Value = 6.983449
$ws.Range("I3").Value = 0.02128501190197005
$ws.Range("J3").Value = 0.02128501190197004
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4472873333333333
$ws.Range("N3").Value = 1.341862
$ws.Range("O3").Value = 0.8927629530472828
$ws.Range("P3").Value = 0.8927629530472827
$ws.Range("Q3").Value = 1.041202760226444
$ws.Range("R3").Value = 9.370824842038
$ws.Range("S3").Value = 0.01900247008124934
$ws.Range("T3").Value = 0.01900247008124933

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ptn"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 103.0385286666667
$ws.Range("H4").Value = 309.115586
$ws.Range("I4").Value = 0.9421603747796319
$ws.Range("J4").Value = 0.9421603747796318
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05372733333333333
$ws.Range("N4").Value = 0.161182
$ws.Range("O4").Value = 0.1072370469527173
$ws.Range("P4").Value = 0.1072370469527173
$ws.Range("Q4").Value = 5.535985375850221
$ws.Range("R4").Value = 49.823868382652
$ws.Range("S4").Value = 0.1010344963472331
$ws.Range("T4").Value = 0.1010344963472331

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ptn"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 103.0385286666667
$ws.Range("H5").Value = 309.115586
$ws.Range("I5").Value = 0.9421603747796319
$ws.Range("J5").Value = 0.9421603747796318
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4472873333333333
$ws.Range("N5").Value = 1.341862
$ws.Range("O5").Value = 0.8927629530472828
$ws.Range("P5").Value = 0.8927629530472827
$ws.Range("Q5").Value = 46.08782871790356
$ws.Range("R5").Value = 414.7904584611319
$ws.Range("S5").Value = 0.8411258784323988
$ws.Range("T5").Value = 0.8411258784323987

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ptn"
$ws.Range("C6").Value = "Ptprz1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.975769333333333
$ws.Range("H6").Value = 11.927308
$ws.Range("I6").Value = 0.03635351138648862
$ws.Range("J6").Value = 0.03635351138648861
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05372733333333333
$ws.Range("N6").Value = 0.161182
$ws.Range("O6").Value = 0.1072370469527173
$ws.Range("P6").Value = 0.1072370469527173
$ws.Range("Q6").Value = 0.2136074842284444
$ws.Range("R6").Value = 1.922467358056
$ws.Range("S6").Value = 0.003898443207449022
$ws.Range("T6").Value = 0.003898443207449022

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ptn"
$ws.Range("C7").Value = "Ptprz1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.975769333333333
$ws.Range("H7").Value = 11.927308
$ws.Range("I7").Value = 0.03635351138648862
$ws.Range("J7").Value = 0.03635351138648861
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4472873333333333
$ws.Range("N7").Value = 1.341862
$ws.Range("O7").Value = 0.8927629530472828
$ws.Range("P7").Value = 0.8927629530472827
$ws.Range("Q7").Value = 1.778311263055111
$ws.Range("R7").Value = 16.004801367496
$ws.Range("S7").Value = 0.0324550681790396
$ws.Range("T7").Value = 0.03245506817903959

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Ptn"
$ws.Range("C8").Value = "Ptprz1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.02199333333333333
$ws.Range("H8").Value = 0.06598
$ws.Range("I8").Value = 0.0002011019319095741
$ws.Range("J8").Value = 0.0002011019319095741
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05372733333333333
$ws.Range("N8").Value = 0.161182
$ws.Range("O8").Value = 0.1072370469527173
$ws.Range("P8").Value = 0.1072370469527173
$ws.Range("Q8").Value = 0.001181643151111111
$ws.Range("R8").Value = 0.01063478836
$ws.Range("S8").Value = 0.00002156557731446916
$ws.Range("T8").Value = 0.00002156557731446916

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Ptn"
$ws.Range("C9").Value = "Ptprz1"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.02199333333333333
$ws.Range("H9").Value = 0.06598
$ws.Range("I9").Value = 0.0002011019319095741
$ws.Range("J9").Value = 0.0002011019319095741
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4472873333333333
$ws.Range("N9").Value = 1.341862
$ws.Range("O9").Value = 0.8927629530472828
$ws.Range("P9").Value = 0.8927629530472827
$ws.Range("Q9").Value = 0.009837339417777778
$ws.Range("R9").Value = 0.08853605475999998
$ws.Range("S9").Value = 0.000179536354595105
$ws.Range("T9").Value = 0.000179536354595105
